# Updates the "cryptos" list (Price / Volume(1h) columns, plus two coin
# rows that swapped rank position) to the values captured in the latest
# GitHub Actions scrape run.
#
# Cells whose new text looks like a plain number (e.g. "208.04") are
# forced to keep their original General/inlineStr text semantics by
# temporarily switching NumberFormat to "@" (Text) before the assignment
# and restoring the prior Style afterward - otherwise Excel's COM layer
# auto-converts the literal into a numeric cell and we lose the original
# formatting/precision (e.g. "208.04" -> 208.0399999999...).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.287.42"
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("D3").Value = "1.577.74"
$ws.Range("E3").Value = "  -0.89%  "
$ws.Range("E4").Value = "  -0.33%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.04"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  -0.35%  "
$ws.Range("E6").Value = "  -1.91%  "
$__style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.31"
$ws.Range("D8").Style = $__style
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("E9").Value = "  -1.53%  "
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "1.801.36"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("D13").Value = "1.581.64"
$ws.Range("E13").Value = "  -0.27%  "
$ws.Range("E14").Value = "  -1.20%  "
$__style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.521"
$ws.Range("D15").Style = $__style
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "27.293.80"
$ws.Range("E16").Value = "  -1.45%  "
$ws.Range("E17").Value = "  -1.03%  "
$__style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.60"
$ws.Range("D18").Style = $__style
$ws.Range("E18").Value = "  -0.96%  "
$__style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("D19").Style = $__style
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  -0.38%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.43"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  -3.55%  "
$ws.Range("E24").Value = "  +1.16%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.90"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  -1.30%  "
$__style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.70"
$ws.Range("D26").Style = $__style
$ws.Range("E26").Value = "  -3.97%  "
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("E30").Value = "  -1.72%  "
$__style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0465"
$ws.Range("D31").Style = $__style
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "1.411.13"
$ws.Range("E33").Value = "  +2.05%  "
$ws.Range("E34").Value = "  -1.46%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +1.54%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.29"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  -1.59%  "
$__style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.939"
$ws.Range("D37").Style = $__style
$ws.Range("E37").Value = "  -2.80%  "
$__style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0166"
$ws.Range("D38").Style = $__style
$ws.Range("E38").Value = "  -2.03%  "
$__style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.823"
$ws.Range("D39").Style = $__style
$ws.Range("E39").Value = "  -0.42%  "
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("E43").Value = "  +3.59%  "
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.35"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.95"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.18"
$ws.Range("D46").Style = $__style
$ws.Range("E46").Value = "  +0.29%  "
$ws.Range("D47").Value = "1.714.12"
$ws.Range("E47").Value = "  -0.93%  "
$__style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.35"
$ws.Range("D48").Style = $__style
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").Value = "0.0₇0990"
$ws.Range("E49").Value = "  -1.25%  "
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0954"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  -1.24%  "
$ws.Range("E51").Value = "  -0.40%  "
